# Update gh-pages output ("杭州-漫展信息.xlsx") to match the newly generated data.
# Sheets (in tab order): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------------------
# 1) 展览 ("展览") — refreshed "想去人数" (F column) counts for existing rows
# ---------------------------------------------------------------------------
$ws1.Cells.Item(4,6).Value = 586
$ws1.Cells.Item(7,6).Value = 1928
$ws1.Cells.Item(8,6).Value = 5330
$ws1.Cells.Item(9,6).Value = 1496
$ws1.Cells.Item(11,6).Value = 3088
$ws1.Cells.Item(13,6).Value = 36
$ws1.Cells.Item(14,6).Value = 1281
$ws1.Cells.Item(15,6).Value = 4229
$ws1.Cells.Item(16,6).Value = 1012
$ws1.Cells.Item(18,6).Value = 1657
$ws1.Cells.Item(19,6).Value = 2598
$ws1.Cells.Item(22,6).Value = 131
$ws1.Cells.Item(23,6).Value = 145
$ws1.Cells.Item(24,6).Value = 970
$ws1.Cells.Item(25,6).Value = 288
$ws1.Cells.Item(27,6).Value = 76
$ws1.Cells.Item(29,6).Value = 1080
$ws1.Cells.Item(30,6).Value = 371
$ws1.Cells.Item(31,6).Value = 38
$ws1.Cells.Item(32,6).Value = 126
$ws1.Cells.Item(34,6).Value = 256
$ws1.Cells.Item(35,6).Value = 1644
$ws1.Cells.Item(36,6).Value = 2165
$ws1.Cells.Item(38,6).Value = 36
$ws1.Cells.Item(39,6).Value = 246
$ws1.Cells.Item(40,6).Value = 603
$ws1.Cells.Item(41,6).Value = 274

# ---------------------------------------------------------------------------
# 2) 展览 — a new event was scraped in ahead of row 44 (2024-07-20, 第五届华盟
#    次元嘉年华), so insert a blank row at 44 and shift the old rows 44-47
#    down to 45-48. Also bump the F value of the two rows right after the
#    inserted one (they were re-scraped too).
# ---------------------------------------------------------------------------
$ws1.Rows.Item(44).Insert()

# Re-apply the bordered/bold/centered index-column formatting (style index 1
# in the original file) to the newly inserted A44 cell.
$idxCell = $ws1.Cells.Item(44,1)
$idxCell.Font.Bold = $true
$idxCell.HorizontalAlignment = -4108   # xlCenter
$idxCell.VerticalAlignment = -4160     # xlTop
$idxCell.Borders.LineStyle = 1

$idxCell.Value = 43
$ws1.Cells.Item(44,2).NumberFormat = "@"
$ws1.Cells.Item(44,2).Value = "2024-07-20"
$ws1.Cells.Item(44,3).Value = "杭州·第五届华盟次元嘉年华&周年庆狂欢"
$ws1.Cells.Item(44,4).Value = "创意路1号 中国智谷富春园区"
$ws1.Cells.Item(44,5).Value = "2024.07.20 10:00-07.21 17:00"
$ws1.Cells.Item(44,6).Value = 0
$ws1.Cells.Item(44,7).Value = 58
$ws1.Cells.Item(44,8).Value = "https://show.bilibili.com/platform/detail.html?id=84762"
$ws1.Cells.Item(44,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/uE6OVg6T1713885553204.jpeg"

# Rows 45 (ex-44, 夏之誓国乙only) and 46 (ex-45, 火影忍者only) were re-scraped
# along with the insert; only their F counts changed.
$ws1.Cells.Item(45,6).Value = 396
$ws1.Cells.Item(46,6).Value = 324

# ---------------------------------------------------------------------------
# 3) 演出 ("演出") — refreshed F counts
# ---------------------------------------------------------------------------
$ws2.Cells.Item(10,6).Value = 144
$ws2.Cells.Item(16,6).Value = 2
$ws2.Cells.Item(18,6).Value = 27

# ---------------------------------------------------------------------------
# 4) 本地生活 ("本地生活") — refreshed F count
# ---------------------------------------------------------------------------
$ws3.Cells.Item(2,6).Value = 734

# ---------------------------------------------------------------------------
# 5) 全部类型 ("全部类型") — refreshed F counts (this sheet is a static merge
#    of the other three and does not get the new row inserted into it)
# ---------------------------------------------------------------------------
$ws4.Cells.Item(2,6).Value = 734
$ws4.Cells.Item(4,6).Value = 586
$ws4.Cells.Item(6,6).Value = 1928
$ws4.Cells.Item(7,6).Value = 5330
$ws4.Cells.Item(8,6).Value = 1496
$ws4.Cells.Item(11,6).Value = 3088
$ws4.Cells.Item(12,6).Value = 36
$ws4.Cells.Item(13,6).Value = 1281
$ws4.Cells.Item(14,6).Value = 4229
$ws4.Cells.Item(15,6).Value = 1012
$ws4.Cells.Item(16,6).Value = 1657
$ws4.Cells.Item(18,6).Value = 2598
$ws4.Cells.Item(24,6).Value = 145
$ws4.Cells.Item(25,6).Value = 144
$ws4.Cells.Item(26,6).Value = 970
$ws4.Cells.Item(27,6).Value = 288
$ws4.Cells.Item(29,6).Value = 76
$ws4.Cells.Item(32,6).Value = 1080
$ws4.Cells.Item(33,6).Value = 371
$ws4.Cells.Item(34,6).Value = 38
$ws4.Cells.Item(35,6).Value = 1644
$ws4.Cells.Item(36,6).Value = 2165
$ws4.Cells.Item(38,6).Value = 36
$ws4.Cells.Item(40,6).Value = 2
$ws4.Cells.Item(41,6).Value = 246
$ws4.Cells.Item(42,6).Value = 603
$ws4.Cells.Item(43,6).Value = 274
$ws4.Cells.Item(45,6).Value = 396
$ws4.Cells.Item(46,6).Value = 324
$ws4.Cells.Item(49,6).Value = 27
